# Reduce disturbance in PLBVF data: shrink the voltage deviations in
# column C (the post-disturbance / recovered voltage profile) by 10x.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PLBVF")
$ws.Activate()

$ws.Range("C3").Value  = 1.001
$ws.Range("C4").Value  = 1.002
$ws.Range("C5").Value  = 1.003
$ws.Range("C6").Value  = 1.002
$ws.Range("C7").Value  = 1.001
$ws.Range("C9").Value  = 0.999
$ws.Range("C10").Value = 0.998
$ws.Range("C11").Value = 0.997
$ws.Range("C12").Value = 0.998
$ws.Range("C13").Value = 0.999

# Leave the selection where the author last left it when saving.
$ws.Range("C14").Select()
